$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'307.02"
$ws.Range("E2").Value = "'1.83%"
$ws.Range("G2").Value = "'14"
$ws.Range("D3").Value = "'35.88"
$ws.Range("E3").Value = "'0.96%"
$ws.Range("G3").Value = "'14"
$ws.Range("D4").Value = "'5.066"
$ws.Range("E4").Value = "'0.26%"
$ws.Range("G4").Value = "'14"
$ws.Range("D5").Value = "'0.08113"
$ws.Range("E5").Value = "'1.49%"
$ws.Range("G5").Value = "'14"
$ws.Range("D6").Value = "'1.942"
$ws.Range("E6").Value = "'0.71%"
$ws.Range("G6").Value = "'14"
$ws.Range("D7").Value = "'4.132"
$ws.Range("E7").Value = "'1.71%"
$ws.Range("G7").Value = "'14"
$ws.Range("D8").Value = "'7.777"
$ws.Range("E8").Value = "'-0.50%"
$ws.Range("G8").Value = "'14"
$ws.Range("D9").Value = "'0.9304"
$ws.Range("E9").Value = "'0.67%"
$ws.Range("G9").Value = "'14"
$ws.Range("D10").Value = "'0.1336"
$ws.Range("E10").Value = "'-4.32%"
$ws.Range("G10").Value = "'14"
$ws.Range("D11").Value = "'0.1914"
$ws.Range("E11").Value = "'1.23%"
$ws.Range("G11").Value = "'14"
$ws.Range("D12").Value = "'0.09230"
$ws.Range("E12").Value = "'-0.15%"
$ws.Range("G12").Value = "'14"
$ws.Range("D13").Value = "'0.03496"
$ws.Range("E13").Value = "'3.40%"
$ws.Range("G13").Value = "'14"
$ws.Range("D14").Value = "'0.09878"
$ws.Range("E14").Value = "'0.06%"
$ws.Range("G14").Value = "'14"
$ws.Range("D15").Value = "'0.001412"
$ws.Range("E15").Value = "'1.31%"
$ws.Range("G15").Value = "'14"
$ws.Range("D16").Value = "'0.005807"
$ws.Range("E16").Value = "'-0.73%"
$ws.Range("G16").Value = "'14"
$ws.Range("D17").Value = "'3.601"
$ws.Range("E17").Value = "'2.75%"
$ws.Range("G17").Value = "'14"
$ws.Range("D18").Value = "'2.928"
$ws.Range("E18").Value = "'-1.07%"
$ws.Range("G18").Value = "'14"
$ws.Range("E19").Value = "'1.25%"
$ws.Range("G19").Value = "'14"
$ws.Range("D20").Value = "'0.1331"
$ws.Range("E20").Value = "'2.20%"
$ws.Range("G20").Value = "'14"
$ws.Range("D21").Value = "'5.237"
$ws.Range("E21").Value = "'3.84%"
$ws.Range("G21").Value = "'14"
$ws.Range("D22").Value = "'0.2594"
$ws.Range("E22").Value = "'7.85%"
$ws.Range("G22").Value = "'14"
$ws.Range("D23").Value = "'0.04399"
$ws.Range("E23").Value = "'-2.05%"
$ws.Range("G23").Value = "'14"
$ws.Range("D24").Value = "'0.001220"
$ws.Range("E24").Value = "'0.42%"
$ws.Range("G24").Value = "'14"
$ws.Range("D25").Value = "'0.004777"
$ws.Range("E25").Value = "'-0.31%"
$ws.Range("G25").Value = "'14"
$ws.Range("E26").Value = "'31.91%"
$ws.Range("G26").Value = "'14"
$ws.Range("D27").Value = "'0.0003122"
$ws.Range("E27").Value = "'3.85%"
$ws.Range("G27").Value = "'14"
$ws.Range("G28").Value = "'14"
$ws.Range("G29").Value = "'14"
$ws.Range("G30").Value = "'14"
$ws.Range("G31").Value = "'14"
$ws.Range("G32").Value = "'14"
$ws.Range("G33").Value = "'14"
$ws.Range("G34").Value = "'14"
$ws.Range("G35").Value = "'14"
$ws.Range("G36").Value = "'14"
$ws.Range("G37").Value = "'14"
$ws.Range("G38").Value = "'14"
$ws.Range("D39").Value = "'0.01995"
$ws.Range("E39").Value = "'4.20%"
$ws.Range("G39").Value = "'14"
$ws.Range("D40").Value = "'0.05047"
$ws.Range("E40").Value = "'6.37%"
$ws.Range("G40").Value = "'14"
$ws.Range("E41").Value = "'16.25%"
$ws.Range("G41").Value = "'14"
$ws.Range("D42").Value = "'0.007615"
$ws.Range("E42").Value = "'3.42%"
$ws.Range("G42").Value = "'14"
$ws.Range("D43").Value = "'0.1382"
$ws.Range("E43").Value = "'3.79%"
$ws.Range("G43").Value = "'14"
$ws.Range("D44").Value = "'0.002095"
$ws.Range("E44").Value = "'-0.80%"
$ws.Range("G44").Value = "'14"
$ws.Range("D45").Value = "'0.01128"
$ws.Range("E45").Value = "'6.91%"
$ws.Range("G45").Value = "'14"
$ws.Range("D46").Value = "'0.00006401"
$ws.Range("E46").Value = "'2.18%"
$ws.Range("G46").Value = "'14"
$ws.Range("E47").Value = "'-0.42%"
$ws.Range("G47").Value = "'14"
$ws.Range("D48").Value = "'64.96"
$ws.Range("E48").Value = "'0.75%"
$ws.Range("G48").Value = "'14"
$ws.Range("D49").Value = "'0.001188"
$ws.Range("E49").Value = "'-28.53%"
$ws.Range("G49").Value = "'14"
$ws.Range("D50").Value = "'0.00002096"
$ws.Range("E50").Value = "'-0.42%"
$ws.Range("G50").Value = "'14"
$ws.Range("D51").Value = "'0.0001996"
$ws.Range("E51").Value = "'-0.42%"
$ws.Range("G51").Value = "'14"
